# Update Mappings 22 Ontologies
# Adds a new "VIMMP_DEF" column (F) to the mapping sheet, populated
# with an empty-list placeholder value "[]" for each existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the styling of the other header cells (E1).
$ws.Range("F1").Value = "VIMMP_DEF"

# New data values for the existing rows.
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"

# Copy formatting (font/border/alignment) from the neighboring header
# cell so the new header cell matches the existing header row style.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
